# Commit checkList and checkItem
#
# Adds two new worksheets ("CheckList" and "Stickers") to the workbook,
# placed immediately after the existing "Attachment" sheet, following the
# same "No. | TESTCASE | METHOD | ... | STATUS_CODE | STATUS_MESSAGE | No. | No."
# layout used by the other test-data sheets in this workbook.

$wb = $excel.ActiveWorkbook
$attachment = $wb.Worksheets.Item("Attachment")

# ---------------------------------------------------------------------
# Create the two new sheets, in order, right after "Attachment"
# ---------------------------------------------------------------------
$checkList = $wb.Worksheets.Add($null, $attachment)
$checkList.Name = "CheckList"

$stickers = $wb.Worksheets.Add($null, $checkList)
$stickers.Name = "Stickers"

# ---------------------------------------------------------------------
# CheckList sheet: A1:H3
# Columns: No. | TESTCASE | METHOD | name | STATUS_CODE | STATUS_MESSAGE | No. | No.
# ---------------------------------------------------------------------

# Pull header/data formatting (fill + centering) from the Attachment sheet
# so the new sheet matches the existing look (shared style indexes).
$attachment.Range("A1:H2").Copy()
$checkList.Range("A1:H2").PasteSpecial(-4122)
$checkList.Range("A2:H2").Copy()
$checkList.Range("A3:H3").PasteSpecial(-4122)

$checkList.Range("A1").Value = "No."
$checkList.Range("B1").Value = "TESTCASE"
$checkList.Range("C1").Value = "METHOD"
$checkList.Range("D1").Value = "name"
$checkList.Range("E1").Value = "STATUS_CODE"
$checkList.Range("F1").Value = "STATUS_MESSAGE"
$checkList.Range("G1").Value = "No."
$checkList.Range("H1").Value = "No."

$checkList.Range("A2").Value = 1
$checkList.Range("B2").Value = "Create a CheckList successfully"
$checkList.Range("C2").Value = "POST"
$checkList.Range("D2").Value = "This task should do:"
$checkList.Range("E2").Value = 200

# ---------------------------------------------------------------------
# Stickers sheet: A1:J3
# Columns: No. | TESTCASE | METHOD | top | left | zIndex | rotate | image | STATUS_CODE | STATUS_MESSAGE
# ---------------------------------------------------------------------
$attachment.Range("A1:H2").Copy()
$stickers.Range("A1:J2").PasteSpecial(-4122)
$stickers.Range("A2:J2").Copy()
$stickers.Range("A3:J3").PasteSpecial(-4122)

$stickers.Range("A1").Value = "No."
$stickers.Range("B1").Value = "TESTCASE"
$stickers.Range("C1").Value = "METHOD"
$stickers.Range("D1").Value = "top"
$stickers.Range("E1").Value = "left"
$stickers.Range("F1").Value = "zIndex"
$stickers.Range("G1").Value = "rotate"
$stickers.Range("H1").Value = "image"
$stickers.Range("I1").Value = "STATUS_CODE"
$stickers.Range("J1").Value = "STATUS_MESSAGE"

$stickers.Range("A2").Value = 1
$stickers.Range("B2").Value = "Create a Stickers successfully"
$stickers.Range("C2").Value = "POST"
$stickers.Range("D2").Value = 50
$stickers.Range("E2").Value = 50
$stickers.Range("F2").Value = 50
$stickers.Range("G2").Value = 50
$stickers.Range("I2").Value = 200

# Leave "Stickers" as the selected/active tab (matches the authored file).
$stickers.Activate()
